$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1524.8235
$ws.Range("I19").Value = 2319.7144
$ws.Range("J19").Value = 968.4
$ws.Range("K19").Value = 2319.7144
$ws.Range("L19").Value = 968.4
$ws.Range("M19").Value = -2144.7144
$ws.Range("N19").Value = -1318.4

$ws.Range("H33").Value = 47630350
$ws.Range("I33").Value = 71431070
$ws.Range("J33").Value = 28912.572
$ws.Range("K33").Value = 71431070
$ws.Range("L33").Value = 28912.572
$ws.Range("M33").Value = -71430841
$ws.Range("N33").Value = -29370.572

$ws.Range("H69").Value = 3557.1428
$ws.Range("J69").Value = 3557.1428
$ws.Range("L69").Value = 10671.4284
$ws.Range("N69").Value = -12419.4284

$ws.Range("H72").Value = 3557.1428
$ws.Range("J72").Value = 3557.1428
$ws.Range("L72").Value = 32014.2852
$ws.Range("N72").Value = -40750.2852

$ws.Range("H93").Value = 34234
$ws.Range("J93").Value = 34234
$ws.Range("L93").Value = 34234
$ws.Range("N93").Value = -39226

$ws.Range("H101").Value = 308
$ws.Range("I101").Value = 308
$ws.Range("K101").Value = 924
$ws.Range("M101").Value = 698

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5848.5083
$ws.Range("I32").Value = 2984.9592
$ws.Range("J32").Value = 19879.9
$ws.Range("K32").Value = 2984.9592
$ws.Range("L32").Value = 19879.9
$ws.Range("M32").Value = -2697.9592
$ws.Range("N32").Value = -20453.9

$ws.Range("H61").Value = 1405.5238
$ws.Range("I61").Value = 1171.7142
$ws.Range("J61").Value = 1873.1428
$ws.Range("K61").Value = 1171.7142
$ws.Range("L61").Value = 1873.1428
$ws.Range("M61").Value = -959.7141999999999
$ws.Range("N61").Value = -2297.1428

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H122").Value = 2240
$ws.Range("I122").Value = 2164.2144
$ws.Range("J122").Value = 2321.6155
$ws.Range("K122").Value = 6492.6432
$ws.Range("L122").Value = 6964.8465
$ws.Range("M122").Value = -4042.6432
$ws.Range("N122").Value = -11864.8465

$ws.Range("H136").Value = 1405.5238
$ws.Range("I136").Value = 1171.7142
$ws.Range("J136").Value = 1873.1428
$ws.Range("K136").Value = 3515.1426
$ws.Range("L136").Value = 5619.428400000001
$ws.Range("M136").Value = -965.1425999999997
$ws.Range("N136").Value = -10719.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1782.6428
$ws.Range("I86").Value = 1874.625
$ws.Range("J86").Value = 1660
$ws.Range("K86").Value = 1874.625
$ws.Range("L86").Value = 1660
$ws.Range("M86").Value = -751.625
$ws.Range("N86").Value = -3906

$ws.Range("H88").Value = 40000
$ws.Range("J88").Value = 40000
$ws.Range("L88").Value = 40000
$ws.Range("N88").Value = -40812

$ws.Range("H89").Value = 1782.6428
$ws.Range("I89").Value = 1874.625
$ws.Range("J89").Value = 1660
$ws.Range("K89").Value = 9373.125
$ws.Range("L89").Value = 8300
$ws.Range("M89").Value = -3757.125
$ws.Range("N89").Value = -19532

$ws.Range("H91").Value = 40000
$ws.Range("J91").Value = 40000
$ws.Range("L91").Value = 40000
$ws.Range("N91").Value = -42808

$ws.Range("H94").Value = 1611.1111
$ws.Range("I94").Value = 875
$ws.Range("J94").Value = 2200
$ws.Range("K94").Value = 875
$ws.Range("L94").Value = 2200
$ws.Range("M94").Value = -424
$ws.Range("N94").Value = -3102

$ws.Range("H99").Value = 1368.4706
$ws.Range("I99").Value = 940.125
$ws.Range("J99").Value = 1749.2222
$ws.Range("K99").Value = 940.125
$ws.Range("L99").Value = 1749.2222
$ws.Range("M99").Value = 557.875
$ws.Range("N99").Value = -4745.2222

$ws.Range("H107").Value = 1437
$ws.Range("I107").Value = 1437
$ws.Range("K107").Value = 1437
$ws.Range("M107").Value = 483

$ws.Range("H109").Value = 47596
$ws.Range("J109").Value = 47596
$ws.Range("L109").Value = 47596
$ws.Range("N109").Value = -50370

$ws.Range("H134").Value = 15626413
$ws.Range("I134").Value = 19232118
$ws.Range("K134").Value = 57696354
$ws.Range("M134").Value = -57693819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18522132
$ws.Range("I31").Value = 41669468
$ws.Range("K31").Value = 41669468
$ws.Range("M31").Value = -41669173

$ws.Range("H34").Value = 18522132
$ws.Range("I34").Value = 41669468
$ws.Range("K34").Value = 41669468
$ws.Range("M34").Value = -41669266

$ws.Range("H58").Value = 1730.909
$ws.Range("I58").Value = 2512
$ws.Range("J58").Value = 1652.8
$ws.Range("K58").Value = 2512
$ws.Range("L58").Value = 1652.8
$ws.Range("M58").Value = -2309
$ws.Range("N58").Value = -2058.8

$ws.Range("H62").Value = 23883560
$ws.Range("I62").Value = 33435644
$ws.Range("J62").Value = 3350
$ws.Range("K62").Value = 33435644
$ws.Range("L62").Value = 3350
$ws.Range("M62").Value = -33435020
$ws.Range("N62").Value = -4598

$ws.Range("H65").Value = 23883560
$ws.Range("I65").Value = 33435644
$ws.Range("J65").Value = 3350
$ws.Range("K65").Value = 167178220
$ws.Range("L65").Value = 16750
$ws.Range("M65").Value = -167175100
$ws.Range("N65").Value = -22990

$ws.Range("H106").Value = 29487.334
$ws.Range("J106").Value = 29487.334
$ws.Range("L106").Value = 29487.334
$ws.Range("N106").Value = -32011.334

$ws.Range("H107").Value = 2196
$ws.Range("I107").Value = 994.4
$ws.Range("J107").Value = 5200
$ws.Range("K107").Value = 994.4
$ws.Range("L107").Value = 5200
$ws.Range("M107").Value = 925.6
$ws.Range("N107").Value = -9040

$ws.Range("H132").Value = 2393.96
$ws.Range("I132").Value = 1754.8667
$ws.Range("J132").Value = 3352.6
$ws.Range("K132").Value = 5264.6001
$ws.Range("L132").Value = 10057.8
$ws.Range("M132").Value = -2734.6001
$ws.Range("N132").Value = -15117.8

$ws.Range("H136").Value = 1730.909
$ws.Range("I136").Value = 2512
$ws.Range("J136").Value = 1652.8
$ws.Range("K136").Value = 7536
$ws.Range("L136").Value = 4958.4
$ws.Range("M136").Value = -4986
$ws.Range("N136").Value = -10058.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3249.875
$ws.Range("J88").Value = 3249.875
$ws.Range("L88").Value = 9749.625
$ws.Range("N88").Value = -10605.625

$ws.Range("H91").Value = 3249.875
$ws.Range("J91").Value = 3249.875
$ws.Range("L91").Value = 9749.625
$ws.Range("N91").Value = -12713.625

$ws.Range("H107").Value = 668.8125
$ws.Range("I107").Value = 525.75
$ws.Range("J107").Value = 716.5
$ws.Range("K107").Value = 1577.25
$ws.Range("L107").Value = 2149.5
$ws.Range("M107").Value = 342.75
$ws.Range("N107").Value = -5989.5

$ws.Range("H131").Value = 48458.832
$ws.Range("J131").Value = 58066.887
$ws.Range("L131").Value = 174200.661
$ws.Range("N131").Value = -184280.661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1167.5
$ws.Range("I97").Value = 999.9091
$ws.Range("K97").Value = 999.9091
$ws.Range("M97").Value = -503.9091

$ws.Range("H107").Value = 1601.5
$ws.Range("I107").Value = 1594
$ws.Range("J107").Value = 1611.5
$ws.Range("K107").Value = 1594
$ws.Range("L107").Value = 1611.5
$ws.Range("M107").Value = 326
$ws.Range("N107").Value = -5451.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 9572.75
$ws.Range("I68").Value = 17678.857
$ws.Range("J68").Value = 3268
$ws.Range("K68").Value = 17678.857
$ws.Range("L68").Value = 3268
$ws.Range("M68").Value = -16929.857
$ws.Range("N68").Value = -4766

$ws.Range("H71").Value = 9572.75
$ws.Range("I71").Value = 17678.857
$ws.Range("J71").Value = 3268
$ws.Range("K71").Value = 88394.285
$ws.Range("L71").Value = 16340
$ws.Range("M71").Value = -84650.285
$ws.Range("N71").Value = -23828

$ws.Range("H82").Value = 2727.3635
$ws.Range("I82").Value = 3034
$ws.Range("J82").Value = 2612.375
$ws.Range("K82").Value = 3034
$ws.Range("L82").Value = 2612.375
$ws.Range("M82").Value = -2673
$ws.Range("N82").Value = -3334.375

$ws.Range("H85").Value = 2727.3635
$ws.Range("I85").Value = 3034
$ws.Range("J85").Value = 2612.375
$ws.Range("K85").Value = 3034
$ws.Range("L85").Value = 2612.375
$ws.Range("M85").Value = -1786
$ws.Range("N85").Value = -5108.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3202
$ws.Range("I62").Value = 3202
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3202
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2578
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3202
$ws.Range("I65").Value = 3202
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16010
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -12890
$ws.Range("N65").ClearContents()

$ws.Range("H92").Value = 29405.555
$ws.Range("J92").Value = 29405.555
$ws.Range("L92").Value = 29405.555
$ws.Range("N92").Value = -34397.555
